$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (rows 194-264) gets one new weekly record inserted at the
# top (row 194), pushing every existing record down by one row (194-264 ->
# 195-265). Use a real row insert so everything below shifts automatically.
$ws.Rows.Item(194).Insert()

# After the insert, the record that used to be on row 263 now lives on row
# 264 - that's the template for the new row 194 (same market/category/
# variety/price data), except it gets a brand-new date.
$lastCol = 18
for ($col = 1; $col -le $lastCol; $col++) {
    if ($col -ne 4) {
        $ws.Cells.Item(194, $col).Value = $ws.Cells.Item(264, $col).Value2
    }
}
$ws.Cells.Item(194, 4).Value = 45009
